# Update cryptocurrency price/volume figures to reflect the latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '37.502.16'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '2.066.94'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '57.65'
$ws.Range('E8').Value = '  -2.09%  '
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0789'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.82'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '2.372.31'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.11'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.764'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.33'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '2.063.98'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '37.471.36'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.15'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.11'
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('D21').Value = '0.0₃0832'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '227.73'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('E25').Value = '  -3.34%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.01'
$ws.Range('E26').Value = '  +5.40%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '169.38'
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('E32').Value = '  -3.16%  '
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.67'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('E37').Value = '  -3.87%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.30'
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '98.15'
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('D42').Value = '1.490.72'
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0955'
$ws.Range('E43').Value = '  -2.50%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '17.02'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.19'
$ws.Range('E46').Value = '  +3.36%  '
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.03'
$ws.Range('E48').Value = '  -4.17%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.28'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').Value = '2.255.02'
$ws.Range('E51').Value = '  -0.92%  '
